$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.674.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.676.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.64%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.528"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.20"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.51%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0641"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.35%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.915.24"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.674.52"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.602"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.07"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +8.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.00"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.670.65"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.91"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.81"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0720"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.23"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.28%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.17"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.79"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.16%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.68"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0493"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.59%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.55%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.30"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.503.20"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.34%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.32%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.02"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.17%  "

$ws.Range("B37").Value = "Aave"
$ws.Range("C37").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "82.92"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +9.69%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0178"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.32%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.03%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.02"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.838"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.26%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.02"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.26%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.808.18"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "50.06"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "93.14"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0116"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.12%  "
